$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (코칩), shifting all existing
# data rows (2-14) down by one (they become rows 3-15).
$ws.Rows.Item(2).Insert()

# The text-like date columns (A, D, E) must be stored as literal text
# (shared strings), not auto-converted to Excel date serials, so force
# the "Text" number format before writing them, then restore the
# cell style back to Normal so no stray formatting is left behind.
$dateCells = @("A2", "D2", "E2")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New row 2: HD현대마린솔루션 (HD Hyundai Marine Solution) IPO entry.
$ws.Range("A2").Value = "2024-04-25"
$ws.Range("B2").Value = "HD현대마린솔루션"
$ws.Range("C2").Value = "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점"
$ws.Range("D2").Value = "2024-04-30"
$ws.Range("E2").Value = "2024-05-08"
$ws.Range("F2").Value = 742260000
$ws.Range("G2").Value = 8900000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 73300
$ws.Range("J2").Value = 83400
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 83400
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 50
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "255.37 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# The row-insert copies the formatting of the row above (bold header
# style) onto the new row; put the whole row back to the plain,
# unstyled look shared by every other data row.
$ws.Range("A2:T2").Style = "Normal"

# The 청약일 (subscription date) for 코칩, now on row 3, moves from
# 2024-04-24 to 2024-04-25.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2024-04-25"
$ws.Range("A3").Style = "Normal"
